$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 4 new rows after the header row for the new "Model 1".."Model 4"
#    sensitivity results, and drop the 4 "Bayesian Model" rows that are no
#    longer part of the table (they are currently the last 4 data rows, which
#    after the insert live at rows 13:16).
# ---------------------------------------------------------------------------
$ws.Rows("2:5").Insert()
$ws.Rows("13:16").Delete()

# ---------------------------------------------------------------------------
# 2. Rename the two "% Difference" column headers.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,4).Value = "% Difference in Linf fits"
$ws.Cells.Item(1,7).Value = "% Difference in K fits"

# ---------------------------------------------------------------------------
# 3. Fill in the full data block (rows 2-12) with the new sensitivity values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "Model 1"
$ws.Cells.Item(2,2).Value = 53.97
$ws.Cells.Item(2,3).Value = 61.26
$ws.Cells.Item(2,4).Value = -11.91
$ws.Cells.Item(2,5).Value = 0.46
$ws.Cells.Item(2,6).Value = 0.3
$ws.Cells.Item(2,7).Value = 56.61

$ws.Cells.Item(3,1).Value = "Model 2"
$ws.Cells.Item(3,2).Value = 54.08
$ws.Cells.Item(3,3).Value = 61.79
$ws.Cells.Item(3,4).Value = -12.48
$ws.Cells.Item(3,5).Value = 0.46
$ws.Cells.Item(3,6).Value = 0.29
$ws.Cells.Item(3,7).Value = 59.43

$ws.Cells.Item(4,1).Value = "Model 3"
$ws.Cells.Item(4,2).Value = 56.6
$ws.Cells.Item(4,3).Value = 73.69
$ws.Cells.Item(4,4).Value = -23.18
$ws.Cells.Item(4,5).Value = 0.45
$ws.Cells.Item(4,6).Value = 0.17
$ws.Cells.Item(4,7).Value = 158.6

$ws.Cells.Item(5,1).Value = "Model 4"
$ws.Cells.Item(5,2).Value = 56.71
$ws.Cells.Item(5,3).Value = 73.67
$ws.Cells.Item(5,4).Value = -23.03
$ws.Cells.Item(5,5).Value = 0.44
$ws.Cells.Item(5,6).Value = 0.18
$ws.Cells.Item(5,7).Value = 151.11

$ws.Cells.Item(6,1).Value = "Model 5"
$ws.Cells.Item(6,2).Value = 54.5
$ws.Cells.Item(6,3).Value = 62.95
$ws.Cells.Item(6,4).Value = -13.43
$ws.Cells.Item(6,5).Value = 0.45
$ws.Cells.Item(6,6).Value = 0.27
$ws.Cells.Item(6,7).Value = 62.07

$ws.Cells.Item(7,1).Value = "Model 6"
$ws.Cells.Item(7,2).Value = 69.92
$ws.Cells.Item(7,3).Value = 77.96
$ws.Cells.Item(7,4).Value = -10.32
$ws.Cells.Item(7,5).Value = 0.2
$ws.Cells.Item(7,6).Value = 0.12
$ws.Cells.Item(7,7).Value = 63.4

$ws.Cells.Item(8,1).Value = "Model 7"
$ws.Cells.Item(8,2).Value = 64.93
$ws.Cells.Item(8,3).Value = 64.74
$ws.Cells.Item(8,4).Value = 0.29
$ws.Cells.Item(8,5).Value = 0.27
$ws.Cells.Item(8,6).Value = 0.26
$ws.Cells.Item(8,7).Value = 3.57

$ws.Cells.Item(9,1).Value = "Model 8"
$ws.Cells.Item(9,2).Value = 54.5
$ws.Cells.Item(9,3).Value = 66.89
$ws.Cells.Item(9,4).Value = -18.53
$ws.Cells.Item(9,5).Value = 0.45
$ws.Cells.Item(9,6).Value = 0.25
$ws.Cells.Item(9,7).Value = 76.05

$ws.Cells.Item(10,1).Value = "Model 9"
$ws.Cells.Item(10,2).Value = 54.5
$ws.Cells.Item(10,3).Value = 64.74
$ws.Cells.Item(10,4).Value = -15.82
$ws.Cells.Item(10,5).Value = 0.45
$ws.Cells.Item(10,6).Value = 0.26
$ws.Cells.Item(10,7).Value = 70.66

$ws.Cells.Item(11,1).Value = "Model 10"
$ws.Cells.Item(11,2).Value = 54.5
$ws.Cells.Item(11,3).Value = 69.34
$ws.Cells.Item(11,4).Value = -21.41
$ws.Cells.Item(11,5).Value = 0.45
$ws.Cells.Item(11,6).Value = 0.15
$ws.Cells.Item(11,7).Value = 205.08

$ws.Cells.Item(12,1).Value = "Model 11"
$ws.Cells.Item(12,2).Value = 66.47
$ws.Cells.Item(12,3).Value = 68.14
$ws.Cells.Item(12,4).Value = -2.45
$ws.Cells.Item(12,5).Value = 0.25
$ws.Cells.Item(12,6).Value = 0.21
$ws.Cells.Item(12,7).Value = 15.47

# ---------------------------------------------------------------------------
# 4. Formatting touch-ups.
# ---------------------------------------------------------------------------
# The "Model" label column is no longer italicised.
$ws.Range("A2:A12").Font.Italic = $false

# The last column ("% Difference in K fits") now shows one decimal place.
$ws.Range("G2:G12").NumberFormat = "0.0"

# Column A is widened to fit the longer "Model" labels.
$ws.Columns.Item(1).ColumnWidth = 26.69

# Match the saved selection from the authored workbook.
$ws.Range("B16").Select()
